# SWS.xlsx update — refresh unit prices / reference codes on rows 12, 13, 55, 56
# across the summary sheet ("BON DE PREPARATION") and every supplier sheet ("1".."19"),
# then restore each sheet's selection and finish with sheet "4" active —
# mirroring the click-through-every-tab session captured in the commit.

$wb = $excel.ActiveWorkbook

# Every sheet in the workbook (summary + the 19 numbered tabs) carries an
# identical copy of these five cells; update them all.
$allSheetNames = @("BON DE PREPARATION","1","2","3","4","5","6","7","8","9","10","11","12","13","14","15","16","17","18","19")

foreach ($name in $allSheetNames) {
    $sh = $wb.Worksheets.Item($name)
    $sh.Range("A12").Value = 2993253
    $sh.Range("D12").Value = 216.5
    $sh.Range("D13").Value = 216.5
    $sh.Range("D55").Value = 98
    $sh.Range("D56").Value = 228
}

# Restore the per-sheet selection/active-cell exactly as captured after the edit pass.
# The summary sheet's cursor ends on B15 after having scrolled back up to the top;
# every numbered sheet's cursor lands on A12 (where the data was just touched),
# except "1" whose prior selection (F62) also collapses to A12.
$ws = $wb.Worksheets.Item("BON DE PREPARATION")
[void]$ws.Activate()
[void]$ws.Range("B15").Select()

$orderedTabs = @("9","10","11","12","13","14","15","16","17","18","1","19","2","3","4","5","6","7","8")
foreach ($name in $orderedTabs) {
    $sh = $wb.Worksheets.Item($name)
    [void]$sh.Activate()
    [void]$sh.Range("A12").Select()
}

# The session finishes with sheet "4" selected/active (workbookView activeTab points at it).
$final = $wb.Worksheets.Item("4")
[void]$final.Activate()
[void]$final.Range("A12").Select()
